$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Texture Rotation API" (row 9) and "Multitexture API" (row 10) rows,
# plus one extra blank spacer row further down (original row 26, between the
# "Input Method API" row and the "Lightning 2.0 Pre-Final Tasks" header).
# Only columns A:B are shifted up (column C holds an unrelated formatted
# placeholder cell that stays put), so we delete ranges rather than whole rows,
# and work from the bottom up to keep the earlier row numbers valid.
$ws.Range("A26:B26").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)
$ws.Range("A9:B10").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)

# "Fix remaining consistency issues" is now on row 9; mark it as completed.
$ws.Range("B9").Value = "Yes"

# Update selection to match the saved view state.
$ws.Range("A11").Select()
